$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305608987808228
$ws.Range("B1").Value = 3.74807071685791
$ws.Range("C1").Value = 3.922101020812988
$ws.Range("D1").Value = 3.004763126373291
$ws.Range("E1").Value = 1.047449350357056
